# Insert a new weekly record as row 598 in the "Feria Lagunitas de Puerto
# Montt - Cebolla" sheet. All rows from 598 downward shift down by one,
# which Rows.Insert() handles automatically (including carrying the date
# number-format style down onto the new row's D cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push existing row 598 (and everything below it) down by one row.
$ws.Rows.Item(598).Insert()

# Populate the newly inserted row 598 with the new weekly record.
$ws.Range("A598").Value = 4
$ws.Range("B598").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C598").Value = "Los Lagos"
$ws.Range("D598").Value = 44783
$ws.Range("E598").Value = 10
$ws.Range("F598").Value = 100112004
$ws.Range("G598").Value = "Cebolla"
$ws.Range("H598").Value = "Sin especificar"
$ws.Range("I598").Value = "1a (guarda)"
$ws.Range("J598").Value = 150
$ws.Range("K598").Value = 9500
$ws.Range("L598").Value = 9500
$ws.Range("M598").Value = 9500
$ws.Range("N598").Value = '$/malla 18 kilos'
$ws.Range("O598").Value = "Región de O'Higgins"
$ws.Range("P598").Value = 528
$ws.Range("Q598").Value = 18
$ws.Range("R598").Value = "Hortaliza"
